# Apply roster update to Sheet1: reorder players, update a few rows,
# and append a new row (Jordan Poole) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Guerschon Yabusele", "PF", "Philadelphia 76ers"),
    @("Andre Drummond", "C", "Philadelphia 76ers"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jordan Poole", "PG,SG", "Washington Wizards")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
